# Motor characterization workbook update:
#  - add "volt"/"current"/"wihtout filter" readings block
#  - add old-code vs new-code speed/rpm comparison table
#  - move + resize the existing scatter chart to sit next to the new table
#  - select the last-entered cell (K13), matching the author's final cursor spot

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New small "volt / current / wihtout filter" block (J2:J4, K2:K3) ---
# Written row-5-headers-first, then this block, so the shared-string table
# ends up in the same insertion order as the authored file.
$ws.Range("J5").Value = "old code speed"
$ws.Range("K5").Value = "old code rpm count"
$ws.Range("L5").Value = "new code speed"
$ws.Range("M5").Value = "new code rpm count"

$ws.Range("J2").Value = "volt"
$ws.Range("K2").Value = 4.42

$ws.Range("J3").Value = "current"
$ws.Range("K3").Value = 0.16

$ws.Range("J4").Value = "wihtout filter"

# --- Old-code / new-code comparison table (I6:M9, I10:I13, J13:K13) ---
$ws.Range("I6").Value = 247
$ws.Range("J6").Value = 4320
$ws.Range("K6").Value = 4185
$ws.Range("L6").Value = 4320
$ws.Range("M6").Value = 4185

$ws.Range("I7").Value = 347
$ws.Range("J7").Value = 5160
$ws.Range("K7").Value = 5054

$ws.Range("I8").Value = 447
$ws.Range("J8").Value = 7380
$ws.Range("K8").Value = 7168

$ws.Range("I9").Value = 547
$ws.Range("J9").Value = 9480
$ws.Range("K9").Value = 9536

$ws.Range("I10").Value = 647
$ws.Range("I11").Value = 747
$ws.Range("I12").Value = 847

$ws.Range("I13").Value = 947
$ws.Range("J13").Value = 15660
$ws.Range("K13").Value = 14753

# --- Widen the new columns to fit their (longer) header text ---
# (inputs chosen so the engine's pixel-quantized save lands as close as
#  possible to the authored bestFit widths of 14.57/18.29/15.57/19.29)
$ws.Columns.Item(10).ColumnWidth = 13.583333333333332
$ws.Columns.Item(11).ColumnWidth = 17.416666666666668
$ws.Columns.Item(12).ColumnWidth = 14.583333333333332
$ws.Columns.Item(13).ColumnWidth = 18.416666666666668

# --- Reposition/resize the existing chart next to the new table ---
# (Left/Top computed against the post-edit column widths above, so the
#  saved anchor lands on col 17/row 8 -> col 28/row 33, same as authored.)
$co = $ws.ChartObjects().Item(1)
$co.Left = 1081.3515625
$co.Top = 125.25
$co.Width = 659.3125
$co.Height = 371.2499212598425

# --- Match the author's final selection ---
$ws.Range("K13").Select()
